$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hydropower plant parameters")
$ws.Activate()

# Clarify that these data sources are optional: precipitation, evaporation and
# forced (prescribed environmental/irrigation) outflow can be left blank to
# default to zero, while solar/wind capacity-factor sheets can be left blank
# when the scenario being run has no solar or wind component.
$ws.Range("B5").Value = 'name of corresponding worksheet in the "data" Excel sheets from which to pull data (precipitation flux in kg/m^2/s; hours in rows, years in columns); leave empty for zeros'
$ws.Range("B6").Value = 'name of corresponding worksheet in the "data" Excel sheets from which to pull data (evaporation flux in kg/m^2/s; hours in rows, years in columns); leave empty for zeros'
$ws.Range("B8").Value = 'name of corresponding worksheet in the "data" Excel sheets from which to pull data (hourly solar CF as fraction/percentage; hours in rows, years in columns); leave empty if scenario has no solar'
$ws.Range("B9").Value = 'name of corresponding worksheet in the "data" Excel sheets from which to pull data (hourly wind CF as fraction/percentage; hours in rows, years in columns); leave empty if scenario has no wind'
$ws.Range("B7").Value = 'name of corresponding worksheet in the "data" Excel sheets from which to pull data on prescribed (environmental/irrigation) outflow (hourly in m^3/s; hours in rows, years in columns); leave empty for zeros'

# Mirror the author leaving the cursor on B1 after editing this sheet.
$ws.Range("B1").Select()
